$d = $word.ActiveDocument

# Helper: return the paragraph Index of the Nth paragraph whose text
# (paragraph mark stripped) equals $text.
function Find-ParagraphByText($doc, $text, $occurrence) {
    $seen = 0
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            $seen++
            if ($seen -eq $occurrence) {
                return $i
            }
        }
    }
    return -1
}

# The document contains the same "command list" twice: once as the
# normal list, and again (duplicated) further down as a note. Each
# copy ends with a "Turn off fan" paragraph; a new "Toggle Light"
# line needs to be appended to both copies.

# --- Handle the *second* (later) "Turn off fan" first, so inserting
#     its new paragraph does not shift the index of the first one.
$secondIdx = Find-ParagraphByText $d "Turn off fan" 2
$secondPara = $d.Paragraphs.Item($secondIdx)
$secondPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item($secondIdx + 1).Range.Text = "Toggle Light"

# --- Now the first "Turn off fan".
$firstIdx = Find-ParagraphByText $d "Turn off fan" 1
$firstPara = $d.Paragraphs.Item($firstIdx)
$firstPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstIdx + 1).Range.Text = "Toggle Light"

# --- Relocate the hidden "_GoBack" bookmark. It currently sits at the
#     very start of the second copy's "Turn off air" paragraph; move
#     it to the (now) final, empty paragraph of the document instead
#     (this is where Word leaves it once the surrounding text is no
#     longer the last thing on the page).
$bookmark = $d.Bookmarks.Item("_GoBack")
$bookmark.Delete()

$lastPara = $d.Paragraphs.Last
$target = $lastPara.Range
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target)
